$wb = $excel.ActiveWorkbook

# Sheet ALC, row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2761.5
$ws.Range("J32").Value = 2498
$ws.Range("L32").Value = 2498
$ws.Range("N32").Value = -3150

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2386
$ws.Range("I40").Value = 2095.5557
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2095.5557
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -1920.5557
$ws.Range("N40").Value = -5350

# Sheet ALC, row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3188.6
$ws.Range("J51").Value = 3371.5293
$ws.Range("L51").Value = 3371.5293
$ws.Range("N51").Value = -4339.5293

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2976.1904
$ws.Range("J112").Value = 3146.8667
$ws.Range("L112").Value = 9440.6001
$ws.Range("N112").Value = -11656.6001

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 23281.525
$ws.Range("J138").Value = 3293.7
$ws.Range("L138").Value = 9881.099999999999
$ws.Range("N138").Value = -20161.1

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7767.846
$ws.Range("I61").Value = 7767.846
$ws.Range("K61").Value = 7767.846
$ws.Range("M61").Value = -7555.846

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1915.04
$ws.Range("I88").Value = 1568.3334
$ws.Range("J88").Value = 2024.5264
$ws.Range("K88").Value = 1568.3334
$ws.Range("L88").Value = 2024.5264
$ws.Range("M88").Value = -1162.3334
$ws.Range("N88").Value = -2836.5264

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1915.04
$ws.Range("I91").Value = 1568.3334
$ws.Range("J91").Value = 2024.5264
$ws.Range("K91").Value = 1568.3334
$ws.Range("L91").Value = 2024.5264
$ws.Range("M91").Value = -164.3334
$ws.Range("N91").Value = -4832.5264

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 811.5714
$ws.Range("I102").Value = 625.3333
$ws.Range("K102").Value = 625.3333
$ws.Range("M102").Value = 996.6667

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1391451
$ws.Range("J132").Value = 3433
$ws.Range("L132").Value = 10299
$ws.Range("N132").Value = -15359

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7767.846
$ws.Range("I136").Value = 7767.846
$ws.Range("K136").Value = 23303.538
$ws.Range("M136").Value = -20753.538

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2070.5881
$ws.Range("I86").Value = 2265.111
$ws.Range("K86").Value = 2265.111
$ws.Range("M86").Value = -1142.111

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2070.5881
$ws.Range("I89").Value = 2265.111
$ws.Range("K89").Value = 11325.555
$ws.Range("M89").Value = -5709.555

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6388.6577
$ws.Range("I134").Value = 3858.6216
$ws.Range("J134").Value = 100000
$ws.Range("K134").Value = 11575.8648
$ws.Range("L134").Value = 300000
$ws.Range("M134").Value = -9040.864799999999
$ws.Range("N134").Value = -305070

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4602.52
$ws.Range("I31").Value = 5144.1113
$ws.Range("K31").Value = 5144.1113
$ws.Range("M31").Value = -4849.1113

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4602.52
$ws.Range("I34").Value = 5144.1113
$ws.Range("K34").Value = 5144.1113
$ws.Range("M34").Value = -4942.1113

# Sheet CRP, row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 8362.4
$ws.Range("I94").Value = 51499.5
$ws.Range("K94").Value = 51499.5
$ws.Range("M94").Value = -51048.5

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15237
$ws.Range("I99").Value = 16185.467
$ws.Range("K99").Value = 16185.467
$ws.Range("M99").Value = -14687.467

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 15237
$ws.Range("I126").Value = 16185.467
$ws.Range("K126").Value = 48556.401
$ws.Range("M126").Value = -46086.401

# Sheet CUL, row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 86.07692
$ws.Range("I2").Value = 40.363636
$ws.Range("J2").Value = 337.5
$ws.Range("K2").Value = 242.181816
$ws.Range("L2").Value = 2025
$ws.Range("M2").Value = -129.181816
$ws.Range("N2").Value = -2251

# Sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2755200.5
$ws.Range("I4").Value = 4094001
$ws.Range("K4").Value = 12282003
$ws.Range("M4").Value = -12281891

# Sheet CUL, row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 5033
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 5033
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 15099
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -15665

# Sheet CUL, row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2000
$ws.Range("J81").Value = 2000
$ws.Range("L81").Value = 6000
$ws.Range("N81").Value = -8246

# Sheet CUL, row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2000
$ws.Range("J84").Value = 2000
$ws.Range("L84").Value = 18000
$ws.Range("N84").Value = -29232

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1197
$ws.Range("I132").Value = 909.5714
$ws.Range("J132").Value = 1700
$ws.Range("K132").Value = 8186.1426
$ws.Range("L132").Value = 15300
$ws.Range("M132").Value = -5656.1426
$ws.Range("N132").Value = -20360

# Sheet CUL, row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 10928.579
$ws.Range("I133").Value = 9309.532999999999
$ws.Range("K133").Value = 27928.599
$ws.Range("M133").Value = -22868.599

# Sheet CUL, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 4656.276
$ws.Range("I139").Value = 1924.3684
$ws.Range("J139").Value = 9846.9
$ws.Range("K139").Value = 5773.1052
$ws.Range("L139").Value = 29540.7
$ws.Range("M139").Value = -633.1052
$ws.Range("N139").Value = -39820.7

# Sheet GSM, row 45
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Sheet GSM, row 109
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 100000
$ws.Range("J109").Value = 100000
$ws.Range("L109").Value = 100000
$ws.Range("N109").Value = -102080

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 148176.58
$ws.Range("I122").Value = 168706.17
$ws.Range("J122").Value = 24999
$ws.Range("K122").Value = 506118.51
$ws.Range("L122").Value = 74997
$ws.Range("M122").Value = -503668.51
$ws.Range("N122").Value = -79897

# Sheet GSM, row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 40469.5
$ws.Range("J123").Value = 40469.5
$ws.Range("L123").Value = 40469.5
$ws.Range("N123").Value = -45369.5

# Sheet LTW, row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3710.111
$ws.Range("I40").Value = 4266.3335
$ws.Range("K40").Value = 4266.3335
$ws.Range("M40").Value = -4130.3335

# Sheet LTW, row 74
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 47138.5
$ws.Range("I74").Value = 48112.332
$ws.Range("K74").Value = 48112.332
$ws.Range("M74").Value = -47114.332

# Sheet LTW, row 77
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 47138.5
$ws.Range("I77").Value = 48112.332
$ws.Range("K77").Value = 144336.996
$ws.Range("M77").Value = -139344.996

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3175.25
$ws.Range("I122").Value = 2904
$ws.Range("J122").Value = 3772
$ws.Range("K122").Value = 8712
$ws.Range("L122").Value = 11316
$ws.Range("M122").Value = -6262
$ws.Range("N122").Value = -16216

# Sheet LTW, row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2648374
$ws.Range("I132").Value = 3868923.8
$ws.Range("K132").Value = 11606771.4
$ws.Range("M132").Value = -11604241.4

# Sheet WVR, row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 55287.832
$ws.Range("J123").Value = 55287.832
$ws.Range("L123").Value = 55287.832
$ws.Range("N123").Value = -65087.832

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1779.1154
$ws.Range("I126").Value = 1654.6086
$ws.Range("K126").Value = 4963.825800000001
$ws.Range("M126").Value = -2493.825800000001
